$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# Append the three new daily observations to the Data sheet
$ws.Range("A441").Value = 45117
$ws.Range("B441").Value = 1811.981
$ws.Range("A442").Value = 45118
$ws.Range("B442").Value = 1775.796
$ws.Range("A443").Value = 45119
$ws.Range("B443").Value = 1820.146

# Match the date-column formatting used by the rest of column A
$ws.Range("A440").Copy()
$ws.Range("A441:A443").PasteSpecial(-4122)

# Update the series metadata on the SeriesInfo sheet, keeping the
# date-looking values as plain text (not auto-converted to date serials)
$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Range("B3").Value = "2023-07-12"
$wsInfo.Range("B3").ClearFormats()

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Range("B4").Value = "2023-07-12"
$wsInfo.Range("B4").ClearFormats()

$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Range("B7").Value = "2023-07-12"
$wsInfo.Range("B7").ClearFormats()

$wsInfo.Range("B14").NumberFormat = "@"
$wsInfo.Range("B14").Value = "2023-07-12 13:01:06-05"
$wsInfo.Range("B14").ClearFormats()
